# Dr. Yi Wang @ Fudan
#
# 1) The "datetimeFigureOut" auto-date field cached on the slide master and
#    on every slide layout shows the day the deck was last saved. Advance it
#    by one day: 1/24/2023 -> 1/25/2023.
# 2) The bottom-right headshot picture on slide 1 is resized slightly
#    smaller (keeping its top-left anchor fixed).

$p = $ppt.ActivePresentation

function Update-DatePlaceholders($shapes, $newText) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Type -eq 14 -and $sh.HasTextFrame) {
            if ($sh.PlaceholderFormat.Type -eq 16) {
                $sh.TextFrame.TextRange.Text = $newText
            }
        }
    }
}

$newDate = "1/25/2023"

# Slide master's Date placeholder.
Update-DatePlaceholders $p.SlideMaster.Shapes $newDate

# Every slide layout's Date placeholder.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholders $layout.Shapes $newDate
}

# Resize the fourth picture on slide 1 (the bottom-right headshot),
# keeping its top-left position (Left/Top) unchanged.
$slide = $p.Slides.Item(1)
$pic = $slide.Shapes.Item(4)
$pic.Width = 1391765 / 914400 * 72
$pic.Height = 1762322 / 914400 * 72
